$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Locate the paragraph that ends with "...robust enough to solve
# these issues." (currently holds the _GoBack bookmark) — this is
# the anchor we insert all the new "12 Mar 2020" class-notes content
# after.
# -----------------------------------------------------------------
$anchorIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.Contains("robust enough to solve these issues")) {
        $anchorIdx = $i
        break
    }
}

$idx = $anchorIdx

# blank paragraph right after the anchor
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"

# "12 Mar 2020" heading
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Heading 3"
$p.Range.Text = "12 Mar 2020"

# blank paragraph
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"

# "Watermarking is a detection problem (yes/no)"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Heading 4"
$p.Range.Text = "Watermarking is a detection problem (yes/no)"

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "While data hiding could deal with hidden data and potential extraction of said data. The slight nuance in definition here. We can consider and think about data-hiding as a communication through an unreliable channel. We have to address capacity."

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Heading 4"
$p.Range.Text = "Watermarking and DH can be prone to security leakages"

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "For instance, there are cases where we can accidentally leak information about the key from the watermark itself, and that is a security consideration that we have to keep in mind."

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Heading 4"
$p.Range.Text = "Stochastics and Perceptual invisibility"

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "Seen from the stego notes, we can skip in favor of the paper notes"

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Heading 4"
$p.Range.Text = "The tradeoffs of security, robustness, and visibility within watermarking, steganography, and data hiding is at the core of the discussion."

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "This is considered an exam question. Explain the tradeoffs, and what each specification of the field is looking for and giving up in order to achieve their goals."

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "We consider:"

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "- Content related/unrelated messages"

# Paragraph with two runs: " " + "- message existence known/unknown"
# Build as two temp paragraphs, then delete the paragraph mark between
# them so they merge into a single paragraph with two separate runs
# (avoids merging the text into one run).
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = " "

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p2 = $d.Paragraphs($idx)
$p2.Style = "Normal"
$p2.Range.Text = "- message existence known/unknown"

$mergePos = $p.Range.End - 1
$d.Range($mergePos, $p.Range.End).Delete()
$idx = $idx - 1
$p = $d.Paragraphs($idx)

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "The applications here on robust watermarking, IP guarantees, authentication, tracking, fragile watermarks, tamperproofing, content recovery, overt embedded communications, media tracking, control of access, etc..."

$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Style = "Normal"
$p.Range.Text = "In a watermarking scheme, we can consider that each transaction will be re-watermarked specifically tailored to that transaction, to state the rights and usage agreements. And we expect that watermark to survive even after attempts at removing it."

# Move the _GoBack bookmark from the old anchor paragraph to the end
# of this, the final newly-added paragraph.
$bmPos = $p.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $p.Range.End))

# -----------------------------------------------------------------
# Trim the trailing empty paragraphs: there were 5 blank paragraphs
# after the (old) bookmark paragraph; only 3 should remain.
# -----------------------------------------------------------------
$lastContentIdx = $idx
$trailing1 = $d.Paragraphs($lastContentIdx + 1)
$trailing1.Range.Delete()
$trailing2 = $d.Paragraphs($lastContentIdx + 1)
$trailing2.Range.Delete()
